$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "74.135.35"
Set-TextValue "E2" "  +6.40%  "
Set-TextValue "D3" "2.642.63"
Set-TextValue "E3" "  +8.07%  "
Set-TextValue "E4" "  +0.00%  "
Set-TextValue "D5" "185.81"
Set-TextValue "E5" "  +11.99%  "
Set-TextValue "D6" "582.56"
Set-TextValue "E6" "  +3.12%  "
Set-TextValue "D7" "0.999"
Set-TextValue "E7" "  -0.07%  "
Set-TextValue "D8" "0.530"
Set-TextValue "E8" "  +3.77%  "
Set-TextValue "D9" "0.191"
Set-TextValue "E9" "  +8.75%  "
Set-TextValue "D10" "2.641.02"
Set-TextValue "E10" "  +8.01%  "
Set-TextValue "E11" "  +0.73%  "
Set-TextValue "D12" "0.353"
Set-TextValue "E12" "  +5.54%  "
Set-TextValue "D13" "4.68"
Set-TextValue "E13" "  +0.15%  "
Set-TextValue "D14" "3.151.83"
Set-TextValue "E14" "  +9.05%  "
Set-TextValue "D15" "74.110.92"
Set-TextValue "E15" "  +6.59%  "
Set-TextValue "D16" "0.0000185"
Set-TextValue "E16" "  +1.59%  "
Set-TextValue "D17" "26.11"
Set-TextValue "E17" "  +8.96%  "
Set-TextValue "D18" "2.653.08"
Set-TextValue "E18" "  +8.76%  "
Set-TextValue "E19" "  +30.20%  "
Set-TextValue "D20" "11.82"
Set-TextValue "E20" "  +9.77%  "
Set-TextValue "D21" "367.53"
Set-TextValue "E21" "  +7.55%  "
Set-TextValue "E22" "  +11.86%  "
Set-TextValue "E23" "  +4.13%  "
Set-TextValue "D24" "6.19"
Set-TextValue "E24" "  +2.21%  "
Set-TextValue "E25" "  -0.07%  "
Set-TextValue "D26" "69.68"
Set-TextValue "E26" "  +5.17%  "
Set-TextValue "D27" "4.06"
Set-TextValue "E27" "  +4.84%  "
Set-TextValue "B28" "Aptos"
Set-TextValue "C28" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D28" "9.23"
Set-TextValue "E28" "  +8.24%  "
Set-TextValue "B29" "WrappedeETH"
Set-TextValue "C29" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D29" "2.783.99"
Set-TextValue "E29" "  +8.27%  "
Set-TextValue "E30" "  -2.14%  "
Set-TextValue "E31" "  +8.34%  "
Set-TextValue "D32" "514.32"
Set-TextValue "E32" "  +14.80%  "
Set-TextValue "E33" "  +10.11%  "
Set-TextValue "D34" "7.59"
Set-TextValue "E34" "  +3.83%  "
Set-TextValue "E35" "  +6.76%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  +0.04%  "
Set-TextValue "D37" "162.75"
Set-TextValue "E37" "  +0.44%  "
Set-TextValue "E38" "  +7.11%  "
Set-TextValue "E39" "  +5.23%  "
Set-TextValue "D40" "19.28"
Set-TextValue "E40" "  +1.16%  "
Set-TextValue "E41" "  +0.04%  "
Set-TextValue "D42" "4.87"
Set-TextValue "E42" "  +9.37%  "
Set-TextValue "B43" "PolygonEcosystemToken"
Set-TextValue "C43" "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue "D43" "0.325"
Set-TextValue "E43" "  +7.02%  "
Set-TextValue "B44" "Aave"
Set-TextValue "C44" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D44" "163.51"
Set-TextValue "E44" "  +24.12%  "
Set-TextValue "E45" "  +6.91%  "
Set-TextValue "E46" "  +8.30%  "
Set-TextValue "D47" "38.95"
Set-TextValue "E47" "  +3.44%  "
Set-TextValue "D48" "2.33"
Set-TextValue "E48" "  +8.78%  "
Set-TextValue "E49" "  +16.49%  "
Set-TextValue "D50" "3.60"
Set-TextValue "E50" "  +6.00%  "
Set-TextValue "E51" "  +7.07%  "
